# The source workbook refreshes several rows that are pulled in via an
# external-workbook link ([1]Sheet1!$A1 .. $A7, surfaced on Sheet1 as
# A2:A8). The linked source file went from reporting 13 to reporting 8
# for the first 7 rows of data (rows 8-16 of the source stay at 13, which
# is why Sheet1 rows 9+ keep their original cached value).
#
# The external workbook itself isn't available to this session (it is a
# cloud-hosted file, and links can only be refreshed against an actually
# open/reachable source workbook), so we record the refreshed values
# directly on the dependent cells -- i.e. exactly what "Keep values only"
# / typing the latest numbers over a now-unreachable link would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updatedRows = 2..8
foreach ($r in $updatedRows) {
    $ws.Cells.Item($r, 1).Value2 = 8
}
